# Weekly update: insert one new price-report row at the top of the
# "Espinaca / Vega Modelo de Temuco" block (row 267), pushing the existing
# rows 267-298 down to 268-299.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 267:298 down one row, leaving a blank row 267 (styles/merges
# of the surrounding rows are carried along automatically, same as
# Excel's native "Insert Sheet Rows").
$ws.Rows.Item(267).Insert()

# Populate the newly inserted row with this week's record.
$ws.Cells.Item(267, 1).Value  = 10
$ws.Cells.Item(267, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(267, 3).Value  = "La Araucanía"
$ws.Cells.Item(267, 4).Value  = 45142
$ws.Cells.Item(267, 5).Value  = 9
$ws.Cells.Item(267, 6).Value  = 100112012
$ws.Cells.Item(267, 7).Value  = "Espinaca"
$ws.Cells.Item(267, 8).Value  = "Sin especificar"
$ws.Cells.Item(267, 9).Value  = "Primera"
$ws.Cells.Item(267, 10).Value = 120
$ws.Cells.Item(267, 11).Value = 8000
$ws.Cells.Item(267, 12).Value = 8000
$ws.Cells.Item(267, 13).Value = 8000
$ws.Cells.Item(267, 14).Value = "`$/docena de atados"
$ws.Cells.Item(267, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(267, 16).Value = 2667
$ws.Cells.Item(267, 17).Value = 3
$ws.Cells.Item(267, 18).Value = "Hortaliza"
